$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- Update row 226: mark CTestStepService task as DONE ---
$ws.Range("F226").Value2 = "DONE"
$ws.Range("H226").Value2 = "Completed - 2026-01-16"

# --- Append new backlog rows 272-289 for test execution work ---
$newRows = @(
  "E13|F13.3|US13.3.2|Implement CComponentTestExecution - main execution interface|8|DONE|CComponentTestExecution.java|933 lines - Auto-save, keyboard shortcuts, ISO compliance - 2026-01-16",
  "E13|F13.3|US13.3.2|Update CPageServiceTestRun with component factory method|1|DONE|CPageServiceTestRun.java|createTestExecutionComponent() added - 2026-01-16",
  "E13|F13.3|US13.3.3|Add single-page execution view to CTestRunInitializerService|2|IN_PROGRESS|CTestRunInitializerService.java|Add second view with setAttributeNone(true) for full-screen execution",
  "E13|F13.3|US13.3.3|Add Execute button to test run detail view|1|TODO|Test run detail|Launch button for execution interface",
  "E13|F13.3|US13.3.4|Integrate screenshot capture functionality|3|TODO|CComponentTestExecution|Browser screenshot API for evidence",
  "E13|F13.3|US13.3.4|Implement file attachment upload in execution|2|TODO|CComponentTestExecution|Drag-drop file upload integration",
  "E13|F13.1|US13.1.4|Playwright: Test Cases CRUD operations|3|TODO|automated_tests/|Create, read, update, delete test cases",
  "E13|F13.2|US13.2.1|Playwright: Test Suites CRUD operations|3|TODO|automated_tests/|Create, read, update, delete test suites",
  "E13|F13.3|US13.3.1|Playwright: Test Session creation workflow|2|TODO|automated_tests/|Create session, link to suite",
  "E13|F13.3|US13.3.2|Playwright: Complete test execution workflow|5|TODO|automated_tests/|Execute test, record results, validate statistics",
  "E13|F13.3|US13.3.5|Playwright: Result recording validation|3|TODO|automated_tests/|Test PASS/FAIL/SKIP/BLOCK recording",
  "E13|F13.5|US13.5.1|Create test metrics dashboard component|5|TODO|CComponentTestMetrics|Charts: pass rate, coverage, trends",
  "E13|F13.5|US13.5.2|Implement test report generation service|5|TODO|CTestReportService|Export to PDF/Excel with charts",
  "E13|F13.5|US13.5.3|Implement test coverage matrix by feature|5|TODO|CTestCoverageService|Link test cases to requirements",
  "E13|F13.1|US13.1.1|Run Spotless formatter on all test files|1|TODO|Maven spotless:apply|Format Java files to coding standards",
  "E13|F13.1|US13.1.1|Run CodeQL security scan|1|TODO|GitHub Actions|Security vulnerability check",
  "E13|F13.1|US13.1.1|Update terminology with ISO/ISTQB UI standards|1|DONE|TESTING_TERMINOLOGY_MAPPING.md|UI component standards added - 2026-01-16",
  "E13|F13.1|US13.1.1|Create code pattern compliance audit|1|DONE|CODE_PATTERN_COMPLIANCE_AUDIT.md|100% pattern compliance verified - 2026-01-16"
)

$startRow = 272
$r = $startRow
foreach ($line in $newRows) {
    $parts = $line.Split("|")
    $ws.Cells.Item($r, 1).Value2 = $parts[0]
    $ws.Cells.Item($r, 2).Value2 = $parts[1]
    $ws.Cells.Item($r, 3).Value2 = $parts[2]
    $ws.Cells.Item($r, 4).Value2 = $parts[3]
    $ws.Cells.Item($r, 5).Value2 = [int]$parts[4]
    $ws.Cells.Item($r, 6).Value2 = $parts[5]
    $ws.Cells.Item($r, 7).Value2 = $parts[6]
    $ws.Cells.Item($r, 8).Value2 = $parts[7]
    $r = $r + 1
}
